$wb = $excel.ActiveWorkbook

$wsM = $wb.Worksheets.Item("ScoreM")
$wsF = $wb.Worksheets.Item("ScoreF")

# ScoreM: update individual score values (workout-2 rep counts)
$wsM.Cells.Item(5, 2).Value = 11
$wsM.Cells.Item(10, 2).Value = 21

# ScoreF: fill in newly-recorded workout-5 results (Minute5 / Second5 / Rep5)
# -> columns O, P, Q for rows 2-10
$data = @(
    @(12, 14, 300),
    @(14, 1, 300),
    @(15, 0, 273),
    @(13, 6, 300),
    @(14, 46, 300),
    @(13, 49, 300),
    @(12, 14, 300),
    @(15, 0, 274),
    @(12, 27, 300)
)
$rows = $data.Length
$cols = 3
$arr = New-Object 'object[,]' $rows,$cols
for ($i = 0; $i -lt $rows; $i++) {
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}
$wsF.Range("O2:Q10").Value = $arr

# Update the selection / active-tab state to match the saved workbook:
# ScoreF's selection moves to S12, and ScoreM becomes the active (selected) sheet/tab.
$wsF.Range("S12").Select()
$wsM.Range("D18").Select()
